$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell while forcing Excel to keep it as
# literal text (so numeric-looking strings such as phone numbers with a
# leading zero, "User Rec Id" numbers, or dd-mm-yyyy date strings, are not
# silently auto-converted into numbers / date serials by Excel's type
# inference). A leading apostrophe is the standard Excel text-entry prefix;
# Excel strips it from the stored value and simply remembers the cell held
# text (this reuses the workbook's existing "quote prefix" cell style
# instead of manufacturing a new number format/style).
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# --- Sheet1 : WF1_ScheduleNotify_Hourly row (ETRS / Schedule Notify data) ---
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1 "F2"  "0325555887"
Set-TextValue $ws1 "N2"  "28-05-2024"
Set-TextValue $ws1 "P2"  "31-05-2024 05:00:00 PM"
Set-TextValue $ws1 "AC2" "28-05-2024"
Set-TextValue $ws1 "AE2" "8083598845"
Set-TextValue $ws1 "AN2" "126172"
Set-TextValue $ws1 "AT2" "2212029859"
Set-TextValue $ws1 "AX2" "3284231907"

# --- Sheet2 ---
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2 "F2"  "0325555887"
Set-TextValue $ws2 "AE2" "8083598845"
Set-TextValue $ws2 "AT2" "2212029859"
Set-TextValue $ws2 "AX2" "3284231907"

# --- Sheet3 ---
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3 "F2"  "0325555887"
Set-TextValue $ws3 "AE2" "8083598845"
Set-TextValue $ws3 "AT2" "2212029859"
Set-TextValue $ws3 "AX2" "3284231907"

# --- Sheet4 ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4 "F2"  "0325555887"
Set-TextValue $ws4 "AE2" "8083598845"
Set-TextValue $ws4 "AT2" "2212029859"
Set-TextValue $ws4 "AX2" "3284231907"
